$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.813.82'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '3.382.17'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '3.379.10'
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '3.951.48'
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('E14').Value = '  +1.74%  '
$ws.Range('E15').Value = '  -2.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').Value = '3.374.58'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '61.815.89'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.556'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('D24').Value = '3.520.95'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0770'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.775'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('E47').Value = '  +3.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').Value = '2.377.57'
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0263'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.50%  '
